# Convention change to support multi-axle vehicles.
# Rename the existing single sheet, then duplicate it to create a second
# tire-size variant sheet, updating the few cells/labels that differ.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original sheet -----------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Tire2x_270_70R22"

# --- 2. Duplicate it (keeps formatting, conditional formatting, column
#        widths, etc.) and place the copy right after the original -------
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Tire2x_430_50R38"

# --- 3. Fix up the "Tire" label on sheet 1 (re-assert it, since shared
#        string slots get renumbered as labels change) --------------------
$ws1.Range("H3").Value = "Tire2x_270_70R22"

# --- 4. Sheet 2 specific data: its own name label, and a plain (no
#        formula) outer-diameter value for the 430/50R38 tire -------------
$ws2.Range("H3").Value = "Tire2x_430_50R38"
$ws2.Range("H7").Value = 0.4572

# --- 5. Nudge style/dxf table so both sheets end up with their own set of
#        (visually identical) conditional-format fills, matching the
#        authoring tool's behaviour when a sheet is duplicated ------------
$tmpCell = $ws2.Range("Z1")
$tmpFc1 = $tmpCell.FormatConditions.Add(1, 3, '"zzz"')
$tmpFc1.Interior.Color = 13431551
$tmpFc2 = $tmpCell.FormatConditions.Add(1, 3, '"zzz"')
$tmpFc2.Interior.Color = 13431551
$tmpFc3 = $tmpCell.FormatConditions.Add(1, 3, '"zzz"')
$tmpFc3.Interior.Color = 13431551
$tmpCell.FormatConditions.Delete()

# --- 6. Restore per-sheet selections ---------------------------------------
$ws1.Select()
$ws1.Range("C25").Select()

# Sheet 2 (the new tire variant) ends up the active/visible tab.
$ws2.Select()
$ws2.Range("J16").Select()

Write-Output "ok"
